$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.621.86'
$ws.Range("E2").Value = '  +0.79%  '

# Row 3
$ws.Range("D3").Value = '3.173.63'
$ws.Range("E3").Value = '  +1.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.57'
$ws.Range("E5").Value = '  +1.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.38'
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").Value = '3.176.91'
$ws.Range("E8").Value = '  +2.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.449'
$ws.Range("E9").Value = '  +2.34%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.21'
$ws.Range("E10").Value = '  -1.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.110'
$ws.Range("E11").Value = '  +1.17%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.399'
$ws.Range("E12").Value = '  +3.42%  '

# Row 13
$ws.Range("D13").Value = '3.737.54'
$ws.Range("E13").Value = '  +2.63%  '

# Row 14
$ws.Range("E14").Value = '  +2.73%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.79'
$ws.Range("E15").Value = '  -3.11%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("E16").Value = '  +1.91%  '

# Row 17
$ws.Range("D17").Value = '58.722.63'
$ws.Range("E17").Value = '  +0.82%  '

# Row 18
$ws.Range("D18").Value = '3.194.63'
$ws.Range("E18").Value = '  +2.78%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.16'
$ws.Range("E19").Value = '  +0.76%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.91'
$ws.Range("E20").Value = '  +0.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.07'
$ws.Range("E21").Value = '  -0.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '357.27'
$ws.Range("E22").Value = '  +5.19%  '

# Row 23
$ws.Range("E23").Value = '  +0.14%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.516'
$ws.Range("E24").Value = '  +2.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.71'
$ws.Range("E25").Value = '  +3.85%  '

# Row 26
$ws.Range("E26").Value = '  +1.42%  '

# Row 27
$ws.Range("D27").Value = '0.0₃0951'
$ws.Range("E27").Value = '  +4.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.08%  '

# Row 29
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.55'
$ws.Range("E29").Value = '  +3.74%  '

# Row 30
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.51'
$ws.Range("E30").Value = '  -1.51%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.91'
$ws.Range("E32").Value = '  +2.26%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.39'
$ws.Range("E33").Value = '  +2.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.21'
$ws.Range("E34").Value = '  +0.78%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.90'
$ws.Range("E35").Value = '  +5.97%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.83'
$ws.Range("E36").Value = '  +2.26%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.27'
$ws.Range("E37").Value = '  +3.64%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.75'
$ws.Range("E38").Value = '  -1.97%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.31'
$ws.Range("E39").Value = '  +0.70%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.67'
$ws.Range("E40").Value = '  +13.33%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0676'
$ws.Range("E41").Value = '  +1.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.710'
$ws.Range("E42").Value = '  +4.11%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.05'
$ws.Range("E43").Value = '  +4.26%  '

# Row 44
$ws.Range("D44").Value = '3.226.75'
$ws.Range("E44").Value = '  +2.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '37.00'
$ws.Range("E45").Value = '  +0.50%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0272'
$ws.Range("E46").Value = '  +5.76%  '

# Row 47
$ws.Range("D47").Value = '2.359.58'
$ws.Range("E47").Value = '  +2.22%  '

# Row 48
$ws.Range("E48").Value = '  +0.13%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.03'
$ws.Range("E49").Value = '  +7.59%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.70'
$ws.Range("E50").Value = '  -0.46%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.09'
$ws.Range("E51").Value = '  +1.61%  '
